# Insert a new weekly record as row 11, pushing all subsequent rows
# (old rows 11-93) down by one (to 12-94).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly data point.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across this entire
# data set, so reuse those values for the new row as well.
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "Femacal de La Calera"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44537
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 100112052
$ws.Range("G11").Value = "Albahaca"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 56
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 4000
$ws.Range("N11").Value = "$/docena de matas"
$ws.Range("O11").Value = "Provincia de Quillota"
$ws.Range("P11").Value = 667
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = "Hortaliza"
